$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - add row 5 for the new handback file
# ----------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsOverview.Range("B5").Value = "e2e\517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsOverview.Range("C5").Value = ".md"
$wsOverview.Range("E5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G5").Value = "2017-02-22 08:26:31"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/b354480c3ab030bef42963314777d4b2788f3626/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md", "", "", "e2e\517ad39b-9fcc-4839-9da2-da4ec69df8e6.md") | Out-Null
$wsOverview.Range("B5").Style = "HyperLink"

# ----------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - add row 5 for the new handback file
# ----------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("A5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D5").Value = "e2e"
$wsZhCn.Range("E5").Value = "ht"
$wsZhCn.Range("F5").Value = "True"
$wsZhCn.Range("G5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2017-02-22 08:26:15"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I5").Value = ""
$wsZhCn.Range("J5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsZhCn.Range("K5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.zh-cn.xlf"
$wsZhCn.Range("L5").Value = "2017-02-22 08:27:13"
$wsZhCn.Range("L5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("M5").Value = ""
$wsZhCn.Range("N5").Value = ""
$wsZhCn.Range("O5").Value = "True"
$wsZhCn.Range("P5").Value = ""
$wsZhCn.Range("Q5").Value = "False"
$wsZhCn.Range("R5").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/b354480c3ab030bef42963314777d4b2788f3626/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md", "", "", "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md") | Out-Null
$wsZhCn.Range("A5").Style = "HyperLink"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("J5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-zhcn/blob/b354480c3ab030bef42963314777d4b2788f3626/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md", "", "", "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md") | Out-Null
$wsZhCn.Range("J5").Style = "HyperLink"

# ----------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - add row 5 for the new handback file
# ----------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("A5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D5").Value = "e2e"
$wsDeDe.Range("E5").Value = "ht"
$wsDeDe.Range("F5").Value = "True"
$wsDeDe.Range("G5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.de-de.xlf"
$wsDeDe.Range("H5").Value = "2017-02-22 08:26:31"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I5").Value = ""
$wsDeDe.Range("J5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md"
$wsDeDe.Range("K5").Value = "517ad39b-9fcc-4839-9da2-da4ec69df8e6.b354480c3ab030bef42963314777d4b2788f3626.de-de.xlf"
$wsDeDe.Range("L5").Value = "2017-02-22 08:27:36"
$wsDeDe.Range("L5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("M5").Value = ""
$wsDeDe.Range("N5").Value = ""
$wsDeDe.Range("O5").Value = "True"
$wsDeDe.Range("P5").Value = ""
$wsDeDe.Range("Q5").Value = "False"
$wsDeDe.Range("R5").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/b354480c3ab030bef42963314777d4b2788f3626/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md", "", "", "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md") | Out-Null
$wsDeDe.Range("A5").Style = "HyperLink"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("J5"), "https://github.com/OpenLocalizationTestOrg/ol-test4-dede/blob/b354480c3ab030bef42963314777d4b2788f3626/e2e/517ad39b-9fcc-4839-9da2-da4ec69df8e6.md", "", "", "517ad39b-9fcc-4839-9da2-da4ec69df8e6.md") | Out-Null
$wsDeDe.Range("J5").Style = "HyperLink"

Write-Output "Report generated for handback."
